{"js": "// Office.js (Word JavaScript API) script implementing the commit:\n// \"Fixed missing server main function and minor spelling errors in\n//  Args files and project plan\"\n//\n// Three textual fixes inside the project-plan body text:\n//  1. \"console based\"  -> \"console-based\"      (hyphenate)\n//  2. \"For each command a server\"  -> \"For each command, a server\"  (add comma)\n//  3. \"Meetings. documentation\" -> \"Meetings, documentation\" (typo fix),\n//     together with relocating the stray \"_GoBack\" bookmark that used to\n//     sit inside \"Will use Trello[BOOKMARK] to manage tasks\" so that it now\n//     sits right after \"Meetings,\" (i.e. \"Meetings,[BOOKMARK] documentation\").\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1. \"console based\" -> \"console-based\"\n// ---------------------------------------------------------------------\nconst consoleHits = body.search(\"console based\", { matchCase: true });\nawait context.sync();\nif (consoleHits.items.length > 0) {\n  consoleHits.items[0].insertText(\"console-based\", \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2. \"For each command a server\" -> \"For each command, a server\"\n// ---------------------------------------------------------------------\nconst commandHits = body.search(\"For each command a server\", { matchCase: true });\nawait context.sync();\nif (commandHits.items.length > 0) {\n  commandHits.items[0].insertText(\"For each command, a server\", \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3a. Remove the \"_GoBack\" bookmark from its old spot (inside the\n//     \"Will use Trello / to manage tasks\" run pair).\n// ---------------------------------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-merge \"Will use Trello\" + \" to manage tasks\" into a single run now\n// that nothing splits them.\nconst trelloHits = body.search(\"Will use Trello to manage tasks\", { matchCase: true });\nawait context.sync();\nif (trelloHits.items.length > 0) {\n  trelloHits.items[0].insertText(\"Will use Trello to manage tasks\", \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3b. Fix \"Meetings. documentation\" -> \"Meetings, documentation\"\n// ---------------------------------------------------------------------\nconst meetingsHits = body.search(\"Meetings. documentation\", { matchCase: true });\nawait context.sync();\nif (meetingsHits.items.length > 0) {\n  meetingsHits.items[0].insertText(\"Meetings, documentation\", \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3c. Re-insert \"_GoBack\" immediately after \"Meetings,\" (collapsed point\n//     between the comma and the following space/\"documentation\").\n// ---------------------------------------------------------------------\nconst meetingsCommaHits = body.search(\"Meetings,\", { matchCase: true });\nawait context.sync();\nif (meetingsCommaHits.items.length > 0) {\n  const afterComma = meetingsCommaHits.items[0].getRange(\"End\");\n  afterComma.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script implementing the commit:\n# \"Fixed missing server main function and minor spelling errors in\n#  Args files and project plan\"\n#\n# Three textual fixes inside the project-plan body text:\n#  1. \"console based\"  -> \"console-based\"      (hyphenate)\n#  2. \"For each command a server\"  -> \"For each command, a server\"  (add comma)\n#  3. \"Meetings. documentation\" -> \"Meetings, documentation\" (typo fix),\n#     together with relocating the stray \"_GoBack\" bookmark that used to\n#     sit inside \"Will use Trello[BOOKMARK] to manage tasks\" so that it now\n#     sits right after \"Meetings,\" (i.e. \"Meetings,[BOOKMARK] documentation\").\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1. \"console based\" -> \"console-based\"\n# ---------------------------------------------------------------------\n$find1 = $d.Content.Find\n$find1.Text = \"console based\"\n$find1.Replacement.Text = \"console-based\"\n$find1.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 2. \"For each command a server\" -> \"For each command, a server\"\n# ---------------------------------------------------------------------\n$find2 = $d.Content.Find\n$find2.Text = \"For each command a server\"\n$find2.Replacement.Text = \"For each command, a server\"\n$find2.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 3a. Remove the \"_GoBack\" bookmark from its old spot (inside the\n#     \"Will use Trello / to manage tasks\" run pair).\n# ---------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Re-merge \"Will use Trello\" + \" to manage tasks\" into a single run now\n# that nothing splits them.\n$find3 = $d.Content.Find\n$find3.Text = \"Will use Trello to manage tasks\"\n$find3.Replacement.Text = \"Will use Trello to manage tasks\"\n$find3.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 3b. Fix \"Meetings. documentation\" -> \"Meetings, documentation\"\n# ---------------------------------------------------------------------\n$find4 = $d.Content.Find\n$find4.Text = \"Meetings. documentation\"\n$find4.Replacement.Text = \"Meetings, documentation\"\n$find4.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 3c. Re-insert \"_GoBack\" immediately after \"Meetings,\" (collapsed point\n#     between the comma and the following space/\"documentation\").\n# ---------------------------------------------------------------------\n$find5 = $d.Content.Find\n$find5.Text = \"Meetings,\"\n$find5.Execute() | Out-Null\nif ($find5.Found) {\n    $bmRange = $find5.Parent\n    $bmRange.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
